# Automatic update of files.
# This script applies a row-content permutation to worksheet "Artfynd" (rows 2-22),
# reassigning each row's species-observation record to its new row position, and
# fixes up the small number of empty "touched" marker cells in columns K / AF that
# travel along with the records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBlock = @'
2|A|n|111896643|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575038.7114136803|R|n|6703416.194821274|AW|s|Philipp Weiss|AX|s|Philipp Weiss
3|A|n|111896653|B|n|89183|D|s|LC|E|n|3215|F|s|Rödgul trumpetsvamp|G|s|Craterellus lutescens|H|s|(Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575075.050630242|R|n|6703403.625642136|AW|s|Philipp Weiss|AX|s|Philipp Weiss
4|A|n|111896633|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575100.4050603262|R|n|6703444.118284944|AW|s|Philipp Weiss|AX|s|Philipp Weiss
5|A|n|111884093|B|n|98535|D|s|LC|E|n|222498|F|s|Blåsippa|G|s|Hepatica nobilis|H|s|Schreb.|P|s|Kopparåsen (Kopparåsen), Gstr|Q|n|575065.9914513066|R|n|6703387.648325931|AW|s|Patric Engfeldt|AX|s|Patric Engfeldt
6|A|n|111883983|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kalkberget (Kalkberget), Gstr|Q|n|575058.3527020445|R|n|6703446.206921679|AW|s|Patric Engfeldt|AX|s|Patric Engfeldt
7|A|n|111896654|B|n|89183|D|s|LC|E|n|3215|F|s|Rödgul trumpetsvamp|G|s|Craterellus lutescens|H|s|(Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575072.6962527435|R|n|6703421.833381963|AW|s|Philipp Weiss|AX|s|Philipp Weiss
8|A|n|111896655|B|n|89183|D|s|LC|E|n|3215|F|s|Rödgul trumpetsvamp|G|s|Craterellus lutescens|H|s|(Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575104.6742508161|R|n|6703428.910891063|AW|s|Philipp Weiss|AX|s|Philipp Weiss
9|A|n|111896640|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575025.3556637274|R|n|6703369.042946251|AW|s|Philipp Weiss|AX|s|Philipp Weiss
10|A|n|111896642|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575014.1091647458|R|n|6703387.066676207|AW|s|Philipp Weiss|AX|s|Philipp Weiss
11|A|n|111896637|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575088.0587098968|R|n|6703396.00058554|AW|s|Philipp Weiss|AX|s|Philipp Weiss
12|A|n|111884133|B|n|88899|D|s|NT|E|n|3286|F|s|Flattoppad klubbsvamp|G|s|Clavariadelphus truncatus|H|s|(Quél.) Donk|P|s|Kalkberget (Kalkberget), Gstr|Q|n|575059.034285416|R|n|6703389.477814267|AW|s|Patric Engfeldt|AX|s|Patric Engfeldt
13|A|n|111896644|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575036.4083237475|R|n|6703431.936489306|AW|s|Philipp Weiss|AX|s|Philipp Weiss
14|A|n|111896690|B|n|90687|D|s|LC|E|n|5964|F|s|Fjällig taggsvamp s.str.|G|s|Sarcodon imbricatus s.str.|H|s|(L.:Fr.) P.Karst.|P|s|Kratte masugn, Gstr|Q|n|575060.2881161601|R|n|6703376.67477417|AW|s|Philipp Weiss|AX|s|Philipp Weiss
15|A|n|111896635|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575037.2974304935|R|n|6703389.027347369|AW|s|Philipp Weiss|AX|s|Philipp Weiss
16|A|n|111896639|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575089.384229039|R|n|6703379.745088123|AW|s|Philipp Weiss|AX|s|Philipp Weiss
17|A|n|111896638|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575087.1320314853|R|n|6703393.020834555|AW|s|Philipp Weiss|AX|s|Philipp Weiss
18|A|n|111896636|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575108.85141061|R|n|6703418.142308297|AW|s|Philipp Weiss|AX|s|Philipp Weiss
19|A|n|111896634|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575048.3395925189|R|n|6703452.413791304|AW|s|Philipp Weiss|AX|s|Philipp Weiss
20|A|n|111896641|B|n|90332|D|s|LC|E|n|4769|F|s|Svavelriska|G|s|Lactarius scrobiculatus|H|s|(Scop.:Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575021.3626164712|R|n|6703370.933926445|AW|s|Philipp Weiss|AX|s|Philipp Weiss
21|A|n|111896652|B|n|89183|D|s|LC|E|n|3215|F|s|Rödgul trumpetsvamp|G|s|Craterellus lutescens|H|s|(Fr.) Fr.|P|s|Kratte masugn, Gstr|Q|n|575066.556649723|R|n|6703455.751857814|AW|s|Philipp Weiss|AX|s|Philipp Weiss
22|A|n|111884471|B|n|88899|D|s|NT|E|n|3286|F|s|Flattoppad klubbsvamp|G|s|Clavariadelphus truncatus|H|s|(Quél.) Donk|P|s|Kalkberget (Kalkberget), Gstr|Q|n|575020.8210917887|R|n|6703397.074168184|AW|s|Patric Engfeldt|AX|s|Patric Engfeldt
'@

$rows = $dataBlock -split "`n"
foreach ($rowLine in $rows) {
    $rowLine = $rowLine.Trim()
    if ($rowLine.Length -eq 0) { continue }
    $parts = $rowLine -split '\|'
    $rowNum = [int]$parts[0]
    $i = 1
    while ($i -lt $parts.Length) {
        $col = $parts[$i]
        $kind = $parts[$i + 1]
        $val = $parts[$i + 2]
        $addr = "$col$rowNum"
        if ($kind -eq 'n') {
            $ws.Range($addr).Value2 = [double]$val
        } else {
            $ws.Range($addr).Value2 = $val
        }
        $i += 3
    }
}

# --- Fix up empty marker cells in columns K and AF (no actual value, just presence) ---

# Column K: after the reshuffle the blank marker cell should be attached to rows 5, 6, 12, 22
# (it was previously on rows 2, 6, 14, 17). Row 6 keeps it; remove it from rows that lost it
# and (re)create it on rows that gained it.
$ws.Range("K2").ClearContents()
$ws.Range("K14").ClearContents()
$ws.Range("K17").ClearContents()

$ws.Range("K5").NumberFormat = "General"
$ws.Range("K12").NumberFormat = "General"
$ws.Range("K22").NumberFormat = "General"

# Column AF: the blank marker cell moves from row 16 to row 14.
$ws.Range("AF16").ClearContents()
$ws.Range("AF14").NumberFormat = "General"
